$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "239.87") are written back as text, matching the
# original inlineStr cell type instead of being coerced to numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.339.62'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '239.87'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').Value = '0.6301'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.07437'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').Value = '0.2892'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').Value = '24.94'
$ws.Range('E10').Value = '  +2.32%  '
$ws.Range('D11').Value = '0.07725'
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').Value = '1.842.72'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '4.969'
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').Value = '0.00001026'
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').Value = '81.67'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').Value = '6.236'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '29.333.70'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').Value = '229.12'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('D20').Value = '12.31'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').Value = '7.363'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = '157.95'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').Value = '8.490'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D27').Value = '17.42'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').Value = '0.06941'
$ws.Range('E28').Value = '  +7.14%  '
$ws.Range('D29').Value = '1.461'
$ws.Range('E29').Value = '  +4.86%  '
$ws.Range('D30').Value = '1.481'
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.043'
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.047'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').Value = '1.824'
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').Value = '0.6983'
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = '0.01843'
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('D38').Value = '2.819'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').Value = '1.235.11'
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('D40').Value = '6.812'
$ws.Range('E40').Value = '  +4.14%  '
$ws.Range('D41').Value = '0.9346'
$ws.Range('E41').Value = '  +2.64%  '
$ws.Range('D42').Value = '0.9992'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = '1.986.20'
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('D44').Value = '100.95'
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').Value = '65.25'
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('D46').Value = '0.00000000120'
$ws.Range('E46').Value = '  +3.47%  '
$ws.Range('D47').Value = '7.017'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').Value = '1.702'
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('D49').Value = '8.907'
$ws.Range('E49').Value = '  -1.49%  '
$ws.Range('D50').Value = '0.1138'
$ws.Range('E50').Value = '  -2.79%  '
$ws.Range('D51').Value = '0.3905'
$ws.Range('E51').Value = '  -1.01%  '

# Restore the default cell style so we do not leave a stray
# number-format style applied to the edited range.
$ws.Range("B2:E51").Style = "Normal"
